$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3732.1304
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 3833.5908
$ws.Range("K17").Value = 4500
$ws.Range("L17").Value = 11500.7724
$ws.Range("M17").Value = -4332
$ws.Range("N17").Value = -11836.7724
$ws.Range("H19").Value = 1131.5
$ws.Range("I19").Value = 1062.6
$ws.Range("J19").Value = 1200.4
$ws.Range("K19").Value = 1062.6
$ws.Range("L19").Value = 1200.4
$ws.Range("M19").Value = -887.5999999999999
$ws.Range("N19").Value = -1550.4
$ws.Range("H28").Value = 1500.6666
$ws.Range("I28").Value = 1500.6666
$ws.Range("K28").Value = 1500.6666
$ws.Range("M28").Value = -1015.6666
$ws.Range("H52").Value = 200
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H62").Value = 2749
$ws.Range("I62").Value = 2749
$ws.Range("K62").Value = 2749
$ws.Range("M62").Value = -2125
$ws.Range("H65").Value = 2749
$ws.Range("I65").Value = 2749
$ws.Range("K65").Value = 13745
$ws.Range("M65").Value = -10625
$ws.Range("H116").Value = 8477.4
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H125").Value = 300001730
$ws.Range("I125").Value = 333335230
$ws.Range("K125").Value = 3000017070
$ws.Range("M125").Value = -3000014610

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 11000
$ws.Range("I45").Value = 11000
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 11000
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -10623
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 3510.1072
$ws.Range("I61").Value = 1533.9166
$ws.Range("J61").Value = 4992.25
$ws.Range("K61").Value = 1533.9166
$ws.Range("L61").Value = 4992.25
$ws.Range("M61").Value = -1321.9166
$ws.Range("N61").Value = -5416.25
$ws.Range("H74").Value = 2154.6155
$ws.Range("I74").Value = 2572.6667
$ws.Range("K74").Value = 2572.6667
$ws.Range("M74").Value = -1698.6667
$ws.Range("H77").Value = 2154.6155
$ws.Range("I77").Value = 2572.6667
$ws.Range("K77").Value = 12863.3335
$ws.Range("M77").Value = -8495.333500000001
$ws.Range("H110").Value = 2385
$ws.Range("I110").Value = 2385
$ws.Range("K110").Value = 2385
$ws.Range("M110").Value = -340
$ws.Range("H122").Value = 973.375
$ws.Range("I122").Value = 973.375
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2920.125
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -470.125
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 3510.1072
$ws.Range("I136").Value = 1533.9166
$ws.Range("J136").Value = 4992.25
$ws.Range("K136").Value = 4601.7498
$ws.Range("L136").Value = 14976.75
$ws.Range("M136").Value = -2051.7498
$ws.Range("N136").Value = -20076.75
$ws.Range("H138").Value = 740343.2
$ws.Range("J138").Value = 740343.2
$ws.Range("L138").Value = 740343.2
$ws.Range("N138").Value = -750623.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 671.1429000000001
$ws.Range("I5").Value = 633.1667
$ws.Range("K5").Value = 633.1667
$ws.Range("M5").Value = -520.1667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 779
$ws.Range("I16").Value = 585
$ws.Range("J16").Value = 1555
$ws.Range("K16").Value = 585
$ws.Range("L16").Value = 1555
$ws.Range("M16").Value = -298
$ws.Range("N16").Value = -2129
$ws.Range("H31").Value = 1876.6666
$ws.Range("I31").Value = 1551.6666
$ws.Range("J31").Value = 2310
$ws.Range("K31").Value = 1551.6666
$ws.Range("L31").Value = 2310
$ws.Range("M31").Value = -1256.6666
$ws.Range("N31").Value = -2900
$ws.Range("H34").Value = 1876.6666
$ws.Range("I34").Value = 1551.6666
$ws.Range("J34").Value = 2310
$ws.Range("K34").Value = 1551.6666
$ws.Range("L34").Value = 2310
$ws.Range("M34").Value = -1349.6666
$ws.Range("N34").Value = -2714
$ws.Range("H42").Value = 20000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H58").Value = 1393.7693
$ws.Range("I58").Value = 1393.7693
$ws.Range("K58").Value = 1393.7693
$ws.Range("M58").Value = -1190.7693
$ws.Range("H113").Value = 779
$ws.Range("I113").Value = 585
$ws.Range("J113").Value = 1555
$ws.Range("K113").Value = 585
$ws.Range("L113").Value = 1555
$ws.Range("M113").Value = 1585
$ws.Range("N113").Value = -5895
$ws.Range("H117").Value = 40112
$ws.Range("I117").Value = 40112
$ws.Range("K117").Value = 40112
$ws.Range("M117").Value = -35523
$ws.Range("H136").Value = 1393.7693
$ws.Range("I136").Value = 1393.7693
$ws.Range("K136").Value = 4181.3079
$ws.Range("M136").Value = -1631.3079

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3433.6667
$ws.Range("I34").Value = 154
$ws.Range("J34").Value = 4370.7144
$ws.Range("K34").Value = 462
$ws.Range("L34").Value = 13112.1432
$ws.Range("M34").Value = -378
$ws.Range("N34").Value = -13280.1432
$ws.Range("H39").Value = 6959.357
$ws.Range("J39").Value = 6959.357
$ws.Range("L39").Value = 20878.071
$ws.Range("N39").Value = -21466.071
$ws.Range("H55").Value = 5187
$ws.Range("I55").Value = 1495
$ws.Range("J55").Value = 5714.4287
$ws.Range("K55").Value = 4485
$ws.Range("L55").Value = 17143.2861
$ws.Range("M55").Value = -4308
$ws.Range("N55").Value = -17497.2861
$ws.Range("H134").Value = 7596.727
$ws.Range("I134").Value = 1484.7142
$ws.Range("K134").Value = 4454.142599999999
$ws.Range("M134").Value = 615.8574000000008
$ws.Range("H140").Value = 6303.231
$ws.Range("I140").Value = 868.125
$ws.Range("J140").Value = 14999.4
$ws.Range("K140").Value = 2604.375
$ws.Range("L140").Value = 44998.2
$ws.Range("M140").Value = 2575.625
$ws.Range("N140").Value = -55358.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 14999
$ws.Range("I58").Value = 14999
$ws.Range("K58").Value = 14999
$ws.Range("M58").Value = -14722
$ws.Range("H102").Value = 1041.4375
$ws.Range("J102").Value = 829
$ws.Range("L102").Value = 829
$ws.Range("N102").Value = -4073
$ws.Range("H103").Value = 46666.668
$ws.Range("J103").Value = 46666.668
$ws.Range("L103").Value = 46666.668
$ws.Range("N103").Value = -49010.668
$ws.Range("H113").Value = 4999.3335
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4999.3335
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4999.3335
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9339.333500000001
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1995.7142
$ws.Range("I132").Value = 1990
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5970
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3440
$ws.Range("N132").Value = -11060

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7263.44
$ws.Range("I7").Value = 2833.3333
$ws.Range("K7").Value = 2833.3333
$ws.Range("M7").Value = -2721.3333
$ws.Range("H46").Value = 3041.6667
$ws.Range("I46").Value = 3041.6667
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3041.6667
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2853.6667
$ws.Range("N46").ClearContents()
$ws.Range("H126").Value = 7263.44
$ws.Range("I126").Value = 2833.3333
$ws.Range("K126").Value = 8499.999899999999
$ws.Range("M126").Value = -6029.999899999999
$ws.Range("H132").Value = 2808.7917
$ws.Range("I132").Value = 2644.85
$ws.Range("J132").Value = 3628.5
$ws.Range("K132").Value = 7934.549999999999
$ws.Range("L132").Value = 10885.5
$ws.Range("M132").Value = -5404.549999999999
$ws.Range("N132").Value = -15945.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 73333.336
$ws.Range("J119").Value = 73333.336
$ws.Range("L119").Value = 73333.336
$ws.Range("N119").Value = -83009.336
